$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from H1 (bold, bordered, centered) onto the new
# header cells I1:J1, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the new columns: I is constant 1, J duplicates column H (IP).
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $ws.Cells.Item($r, 8).Value2
}
